# Deleted partner guard, users now make their profiles:
# add a new "DeletePartnerUser" permission row just above the existing
# "ReadPartnerNotification" row (old row 9) on the PartnerPermissions sheet.
# All following rows shift down by one, which Excel's row-insert handles
# automatically (formulas / relative refs re-point themselves).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PartnerPermissions")
$ws.Activate()

# Insert a fresh row above row 9 - rows 9..22 (old) become 10..23 (new).
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the new permission.
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Brisanje korisnika"
$ws.Range("D9").Value = "DeletePartnerUser"
$ws.Range("E9").Formula = '=CONCATENATE("insert into ",$A$1,"(",$B$2,", ",$C$2,", ",$D$2,") values(N''",B9,"'', ",IF(TRIM(C9)<>"","N''"&C9&"''","null"),", N''",D9,"'');")'
$ws.Range("R9").Formula = '=CONCATENATE(D9," = ",A9,",")'

# Column A holds plain sequential literals (not a formula), so the insert
# doesn't renumber the rows that got pushed down - fix that up by hand.
for ($r = 10; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Re-stamp the E/R formulas on every pushed-down row so each becomes its own
# standalone formula instead of staying part of the old E4:E22/R4:R22 shared
# group - that group now only spans E4:E9/R4:R9. A plain re-assignment keeps
# the cell tied to the shared group, so clear it out first.
for ($r = 10; $r -le 23; $r++) {
    $eFormula = $ws.Cells.Item($r, 5).Formula
    $rFormula = $ws.Cells.Item($r, 18).Formula
    $ws.Cells.Item($r, 5).ClearContents()
    $ws.Cells.Item($r, 18).ClearContents()
    $ws.Cells.Item($r, 5).Formula = $eFormula
    $ws.Cells.Item($r, 18).Formula = $rFormula
}

# Match the source workbook's formatting for this row: B/D keep the table
# style, but E/R carry no explicit cell style (same as several other rows
# further down in this sheet).
$ws.Range("E9").Style = "Normal"
$ws.Range("R9").Style = "Normal"

# Leave the selection where the author last clicked after inserting the row.
$ws.Range("D24").Select()
